# "Continuação dos testes de integração dos controllers"
#
# The checklist worksheet (Plan1 / sheet1) tracks progress of backend
# controller work. Item 12 ("Adicionar paginação nos verbos GET", row 13)
# moves from "Não realizado" to "Em andamento" now that its integration
# tests are underway, and the active-cell selection that was left on the
# sheet moves up one row (D11 -> D10) to reflect where the author was
# working next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the status of the "Adicionar paginação nos verbos GET" checklist
# item (row 13) from "Não realizado" to "Em andamento".
$ws.Range("C13").Value = "Em andamento"

# Move the saved selection from D11 to D10.
$ws.Range("D10").Select()
